$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, "D").Value = "'66.206.08"
$ws.Cells.Item(2, "D").Style = "Normal"
$ws.Cells.Item(2, "E").Value = "'  -0.40%  "
$ws.Cells.Item(2, "E").Style = "Normal"
$ws.Cells.Item(3, "D").Value = "'3.334.87"
$ws.Cells.Item(3, "D").Style = "Normal"
$ws.Cells.Item(3, "E").Value = "'  -0.16%  "
$ws.Cells.Item(3, "E").Style = "Normal"
$ws.Cells.Item(4, "D").Value = "'0.998"
$ws.Cells.Item(4, "D").Style = "Normal"
$ws.Cells.Item(4, "E").Value = "'  -0.35%  "
$ws.Cells.Item(4, "E").Style = "Normal"
$ws.Cells.Item(5, "D").Value = "'583.83"
$ws.Cells.Item(5, "D").Style = "Normal"
$ws.Cells.Item(6, "D").Value = "'185.57"
$ws.Cells.Item(6, "D").Style = "Normal"
$ws.Cells.Item(6, "E").Value = "'  -2.51%  "
$ws.Cells.Item(6, "E").Style = "Normal"
$ws.Cells.Item(7, "D").Value = "'0.999"
$ws.Cells.Item(7, "D").Style = "Normal"
$ws.Cells.Item(7, "E").Value = "'  -0.07%  "
$ws.Cells.Item(7, "E").Style = "Normal"
$ws.Cells.Item(8, "D").Value = "'3.329.96"
$ws.Cells.Item(8, "D").Style = "Normal"
$ws.Cells.Item(8, "E").Value = "'  -0.05%  "
$ws.Cells.Item(8, "E").Style = "Normal"
$ws.Cells.Item(9, "D").Value = "'0.577"
$ws.Cells.Item(9, "D").Style = "Normal"
$ws.Cells.Item(9, "E").Value = "'  -2.30%  "
$ws.Cells.Item(9, "E").Style = "Normal"
$ws.Cells.Item(10, "E").Value = "'  -2.54%  "
$ws.Cells.Item(10, "E").Style = "Normal"
$ws.Cells.Item(11, "D").Value = "'0.582"
$ws.Cells.Item(11, "D").Style = "Normal"
$ws.Cells.Item(11, "E").Value = "'  -1.60%  "
$ws.Cells.Item(11, "E").Style = "Normal"
$ws.Cells.Item(12, "D").Value = "'47.19"
$ws.Cells.Item(12, "D").Style = "Normal"
$ws.Cells.Item(12, "E").Value = "'  -1.84%  "
$ws.Cells.Item(12, "E").Style = "Normal"
$ws.Cells.Item(13, "D").Value = "'0.0000269"
$ws.Cells.Item(13, "D").Style = "Normal"
$ws.Cells.Item(13, "E").Value = "'  -1.83%  "
$ws.Cells.Item(13, "E").Style = "Normal"
$ws.Cells.Item(14, "D").Value = "'679.44"
$ws.Cells.Item(14, "D").Style = "Normal"
$ws.Cells.Item(14, "E").Value = "'  +11.63%  "
$ws.Cells.Item(14, "E").Style = "Normal"
$ws.Cells.Item(15, "D").Value = "'3.866.12"
$ws.Cells.Item(15, "D").Style = "Normal"
$ws.Cells.Item(15, "E").Value = "'  -0.18%  "
$ws.Cells.Item(15, "E").Style = "Normal"
$ws.Cells.Item(16, "D").Value = "'8.49"
$ws.Cells.Item(16, "D").Style = "Normal"
$ws.Cells.Item(16, "E").Value = "'  -2.81%  "
$ws.Cells.Item(16, "E").Style = "Normal"
$ws.Cells.Item(17, "D").Value = "'66.338.40"
$ws.Cells.Item(17, "D").Style = "Normal"
$ws.Cells.Item(17, "E").Value = "'  -0.27%  "
$ws.Cells.Item(17, "E").Style = "Normal"
$ws.Cells.Item(18, "B").Value = "Chainlink"
$ws.Cells.Item(18, "C").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(18, "D").Value = "'17.92"
$ws.Cells.Item(18, "D").Style = "Normal"
$ws.Cells.Item(18, "E").Value = "'  -1.20%  "
$ws.Cells.Item(18, "E").Style = "Normal"
$ws.Cells.Item(19, "B").Value = "TRON"
$ws.Cells.Item(19, "C").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(19, "D").Value = "'0.118"
$ws.Cells.Item(19, "D").Style = "Normal"
$ws.Cells.Item(19, "E").Value = "'  -0.62%  "
$ws.Cells.Item(19, "E").Style = "Normal"
$ws.Cells.Item(20, "D").Value = "'3.329.25"
$ws.Cells.Item(20, "D").Style = "Normal"
$ws.Cells.Item(20, "E").Value = "'  -0.30%  "
$ws.Cells.Item(20, "E").Style = "Normal"
$ws.Cells.Item(21, "D").Value = "'11.10"
$ws.Cells.Item(21, "D").Style = "Normal"
$ws.Cells.Item(21, "E").Value = "'  -1.01%  "
$ws.Cells.Item(21, "E").Style = "Normal"
$ws.Cells.Item(22, "D").Value = "'0.899"
$ws.Cells.Item(22, "D").Style = "Normal"
$ws.Cells.Item(22, "E").Value = "'  -2.10%  "
$ws.Cells.Item(22, "E").Style = "Normal"
$ws.Cells.Item(23, "D").Value = "'17.80"
$ws.Cells.Item(23, "D").Style = "Normal"
$ws.Cells.Item(23, "E").Value = "'  -5.79%  "
$ws.Cells.Item(23, "E").Style = "Normal"
$ws.Cells.Item(24, "D").Value = "'103.50"
$ws.Cells.Item(24, "D").Style = "Normal"
$ws.Cells.Item(24, "E").Value = "'  +2.30%  "
$ws.Cells.Item(24, "E").Style = "Normal"
$ws.Cells.Item(25, "E").Value = "'  -2.86%  "
$ws.Cells.Item(25, "E").Style = "Normal"
$ws.Cells.Item(26, "D").Value = "'3.98"
$ws.Cells.Item(26, "D").Style = "Normal"
$ws.Cells.Item(26, "E").Value = "'  -1.76%  "
$ws.Cells.Item(26, "E").Style = "Normal"
$ws.Cells.Item(27, "E").Value = "'  +0.23%  "
$ws.Cells.Item(27, "E").Style = "Normal"
$ws.Cells.Item(28, "D").Value = "'9.45"
$ws.Cells.Item(28, "D").Style = "Normal"
$ws.Cells.Item(28, "E").Value = "'  -3.68%  "
$ws.Cells.Item(28, "E").Style = "Normal"
$ws.Cells.Item(29, "D").Value = "'32.57"
$ws.Cells.Item(29, "D").Style = "Normal"
$ws.Cells.Item(29, "E").Value = "'  +6.34%  "
$ws.Cells.Item(29, "E").Style = "Normal"
$ws.Cells.Item(30, "D").Value = "'8.51"
$ws.Cells.Item(30, "D").Style = "Normal"
$ws.Cells.Item(30, "E").Value = "'  -2.76%  "
$ws.Cells.Item(30, "E").Style = "Normal"
$ws.Cells.Item(31, "D").Value = "'6.81"
$ws.Cells.Item(31, "D").Style = "Normal"
$ws.Cells.Item(31, "E").Value = "'  -0.69%  "
$ws.Cells.Item(31, "E").Style = "Normal"
$ws.Cells.Item(32, "D").Value = "'610.40"
$ws.Cells.Item(32, "D").Style = "Normal"
$ws.Cells.Item(32, "E").Value = "'  +6.66%  "
$ws.Cells.Item(32, "E").Style = "Normal"
$ws.Cells.Item(33, "E").Value = "'  -3.26%  "
$ws.Cells.Item(33, "E").Style = "Normal"
$ws.Cells.Item(34, "D").Value = "'11.12"
$ws.Cells.Item(34, "D").Style = "Normal"
$ws.Cells.Item(34, "E").Value = "'  -0.53%  "
$ws.Cells.Item(34, "E").Style = "Normal"
$ws.Cells.Item(35, "D").Value = "'3.840.01"
$ws.Cells.Item(35, "D").Style = "Normal"
$ws.Cells.Item(35, "E").Value = "'  +3.20%  "
$ws.Cells.Item(35, "E").Style = "Normal"
$ws.Cells.Item(36, "E").Value = "'  -1.37%  "
$ws.Cells.Item(36, "E").Style = "Normal"
$ws.Cells.Item(37, "E").Value = "'  -0.08%  "
$ws.Cells.Item(37, "E").Style = "Normal"
$ws.Cells.Item(38, "D").Value = "'56.04"
$ws.Cells.Item(38, "D").Style = "Normal"
$ws.Cells.Item(38, "E").Value = "'  -2.34%  "
$ws.Cells.Item(38, "E").Style = "Normal"
$ws.Cells.Item(39, "D").Value = "'0.129"
$ws.Cells.Item(39, "D").Style = "Normal"
$ws.Cells.Item(39, "E").Value = "'  -2.30%  "
$ws.Cells.Item(39, "E").Style = "Normal"
$ws.Cells.Item(40, "D").Value = "'2.69"
$ws.Cells.Item(40, "D").Style = "Normal"
$ws.Cells.Item(40, "E").Value = "'  -1.45%  "
$ws.Cells.Item(40, "E").Style = "Normal"
$ws.Cells.Item(41, "D").Value = "'0.0₃0704"
$ws.Cells.Item(41, "D").Style = "Normal"
$ws.Cells.Item(41, "E").Value = "'  -4.25%  "
$ws.Cells.Item(41, "E").Style = "Normal"
$ws.Cells.Item(42, "D").Value = "'3.20"
$ws.Cells.Item(42, "D").Style = "Normal"
$ws.Cells.Item(42, "E").Value = "'  -3.92%  "
$ws.Cells.Item(42, "E").Style = "Normal"
$ws.Cells.Item(43, "D").Value = "'32.70"
$ws.Cells.Item(43, "D").Style = "Normal"
$ws.Cells.Item(43, "E").Value = "'  -4.58%  "
$ws.Cells.Item(43, "E").Style = "Normal"
$ws.Cells.Item(44, "D").Value = "'3.43"
$ws.Cells.Item(44, "D").Style = "Normal"
$ws.Cells.Item(44, "E").Value = "'  +3.32%  "
$ws.Cells.Item(44, "E").Style = "Normal"
$ws.Cells.Item(45, "D").Value = "'0.337"
$ws.Cells.Item(45, "D").Style = "Normal"
$ws.Cells.Item(45, "E").Value = "'  -2.55%  "
$ws.Cells.Item(45, "E").Style = "Normal"
$ws.Cells.Item(46, "E").Value = "'  -3.25%  "
$ws.Cells.Item(46, "E").Style = "Normal"
$ws.Cells.Item(47, "D").Value = "'3.00"
$ws.Cells.Item(47, "D").Style = "Normal"
$ws.Cells.Item(47, "E").Value = "'  -12.35%  "
$ws.Cells.Item(47, "E").Style = "Normal"
$ws.Cells.Item(48, "E").Value = "'  -1.84%  "
$ws.Cells.Item(48, "E").Style = "Normal"
$ws.Cells.Item(49, "B").Value = "ThetaToken"
$ws.Cells.Item(49, "C").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(49, "D").Value = "'2.56"
$ws.Cells.Item(49, "D").Style = "Normal"
$ws.Cells.Item(49, "E").Value = "'  -2.42%  "
$ws.Cells.Item(49, "E").Style = "Normal"
$ws.Cells.Item(50, "B").Value = "FirstDigitalUSD"
$ws.Cells.Item(50, "C").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(50, "D").Value = "'1.00"
$ws.Cells.Item(50, "D").Style = "Normal"
$ws.Cells.Item(50, "E").Value = "'  +0.26%  "
$ws.Cells.Item(50, "E").Style = "Normal"
$ws.Cells.Item(51, "B").Value = "Mantle"
$ws.Cells.Item(51, "C").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(51, "D").Value = "'1.31"
$ws.Cells.Item(51, "D").Style = "Normal"
$ws.Cells.Item(51, "E").Value = "'  +2.10%  "
$ws.Cells.Item(51, "E").Style = "Normal"
